$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of destination row -> source row (both refer to the row numbers
# as they exist BEFORE this edit is applied). Row 1 is the header and is
# left untouched. Data occupies rows 2-26 (last column AY).
$rowMap = @{
    2  = 8
    3  = 9
    4  = 10
    5  = 11
    6  = 12
    7  = 13
    8  = 14
    9  = 15
    10 = 16
    11 = 17
    12 = 18
    13 = 19
    14 = 20
    15 = 21
    16 = 22
    17 = 2
    18 = 23
    19 = 24
    20 = 25
    21 = 3
    22 = 4
    23 = 5
    24 = 6
    25 = 26
    26 = 7
}

$firstRow = 2
$lastRow = 26
$lastCol = "AY"

# Snapshot every source row's full contents (A:AY) before writing anything,
# since this is a permutation (rows read from and written to overlap).
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapshot[$r] = $ws.Range("A$r`:$lastCol$r").Value()
}

# Several columns hold plain text that looks numeric/date/time
# (Antal "23", Startdatum/Slutdatum "2023-04-23", Starttid/Sluttid
# "00:00"), not real numbers/dates. Force those destination cells to
# Text format first so Excel doesn't auto-convert the re-assigned
# strings into numbers or date/time serial numbers.
$textCols = @("I", "Y", "Z", "AA", "AB")
for ($r = $firstRow; $r -le $lastRow; $r++) {
    foreach ($col in $textCols) {
        $ws.Range("$col$r").NumberFormat = "@"
    }
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $ws.Range("A$destRow`:$lastCol$destRow").Value = $snapshot[$srcRow]
}
